$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Internal Name" (column A) / "Industry" (column C) for the
# row that used to describe the Engineering Firm / Engineering industry,
# renaming it to the new Communication Provider / Communication entry.
$ws.Range("A3").Value = "Communication Provider"
$ws.Range("C3").Value = "Communication"

# Move the active selection from C4 to B4.
$ws.Range("B4").Select()
